$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A5").Value = "Login Page"
$ws.Range("B5").Value = "Login"
$ws.Range("C5").Value = "/api/authenticate"
$ws.Range("D5").Value = "POST"
$ws.Range("E5").Value = '{ username: "admin", password: "admin@internship" }'
$ws.Range("F5").Value = '{ status: "success/failure" }'

$ws.Range("D14").Select()
